# Weekly update: insert a new week's pricing rows for
# "Macroferia Regional de Talca" / Pera (Packham's Triumph) right after
# the existing block header row (row 386), pushing the rest of the
# table (old rows 387-406) down to rows 390-409.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows before row 387 (shifts 387.. down by 3).
$ws.Rows.Item(387).Insert()
$ws.Rows.Item(387).Insert()
$ws.Rows.Item(387).Insert()

# Common (constant) values shared by every data row in this block.
$mercadoId = 5
$mercado   = "Macroferia Regional de Talca"
$region    = "Maule"
$codreg    = 7
$tipo      = "Fruta"
$productoId = 100104
$producto   = "Frutos de pepita"
$categoriaId = 100104005
$categoria   = "Pera"
$unidad      = "$/bandeja 18 kilos granel"
$origen      = "Provincia de Curicó"
$kgUnidad    = 18

function Set-DataRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [double]$PrecioKg
    )

    $ws.Cells.Item($Row, 1).Value = $mercadoId
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $tipo
    $ws.Cells.Item($Row, 7).Value = $productoId
    $ws.Cells.Item($Row, 8).Value = $producto
    $ws.Cells.Item($Row, 9).Value = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad

    # Date column keeps the same numeric-date format used elsewhere (s=2).
    $ws.Cells.Item($Row, 4).NumberFormat = $ws.Cells.Item($Row - 1, 4).NumberFormat
}

# NOTE: use positional arguments (named -Param passing is unreliable here).
Set-DataRow 387 44516 "Packham's Triumph" "Especial" 210 12000 12000 12000 667
Set-DataRow 388 44516 "Packham's Triumph" "Primera"  220 10000 10000 10000 556
Set-DataRow 389 44516 "Packham's Triumph" "Segunda"  100 8000  8000  8000  444

Write-Output "Inserted rows 387-389; new dimension should be A1:T409"
